$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 578. This shifts the existing rows 578-631
# down to 579-632, matching the target diff (which shows every row from
# 578 onward taking on the values of the row above it, and a brand new
# row 632 appearing with what used to be row 631's data).
$ws.Range("A578").EntireRow.Insert()

# Populate the newly inserted row 578 with the new data record.
$ws.Range("A578").Value = 6
$ws.Range("B578").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C578").Value = "Metropolitana"
$ws.Range("D578").Value = 45106
$ws.Range("E578").Value = 13
$ws.Range("F578").Value = 100112043
$ws.Range("G578").Value = "Pepino ensalada"
$ws.Range("H578").Value = "Sin especificar"
$ws.Range("I578").Value = "Primera"
$ws.Range("J578").Value = 930
$ws.Range("K578").Value = 12000
$ws.Range("L578").Value = 13000
$ws.Range("M578").Value = 12398
$ws.Range("N578").Value = "$/caja 60 unidades"
$ws.Range("O578").Value = "Región de Arica y Parinacota"
$ws.Range("P578").Value = 207
$ws.Range("Q578").Value = 60
$ws.Range("R578").Value = "Hortaliza"
